# Developer Guide: Fix spelling error
# - Update cached datetimeFigureOut field text from "8/7/2018" to "9/21/18"
#   on the slide master and every slide layout's Date Placeholder.
# - Fix misspelled "CrearCommand" -> "ClearCommand" on slide 2.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -eq "8/7/2018") {
                    $tr.Text = "9/21/18"
                }
            }
        }
    }
}

# Slide master date placeholder
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every custom (slide) layout's date placeholder
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Fix the spelling error on slide 2: "CrearCommand" -> "ClearCommand"
$slide2 = $p.Slides.Item(2)
for ($shi = 1; $shi -le $slide2.Shapes.Count; $shi++) {
    $sh = $slide2.Shapes.Item($shi)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "CrearCommand") {
                $sh.TextFrame.TextRange.Text = "ClearCommand"
            }
        }
    }
}
